# Update "want to go" counts (column F) on the "展览" and "全部类型" sheets
# to reflect newly generated output, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 12801
$wsExhibition.Range("F7").Value = 37
$wsExhibition.Range("F10").Value = 12716
$wsExhibition.Range("F14").Value = 7628
$wsExhibition.Range("F22").Value = 372

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 12801
$wsAll.Range("F8").Value = 37
$wsAll.Range("F11").Value = 12716
$wsAll.Range("F15").Value = 7628
$wsAll.Range("F24").Value = 372
